# Unmap 1.9 snippet and fix pivot text (#194)
#
# The "Snippets" table had a row mapping Worksheet.getRange (snippet
# "excel-range-areas" / colorAllFormulaCells) that needs to be unmapped,
# i.e. removed from the table entirely. Deleting the worksheet row shifts
# every following row up by one, shrinks the table/dimension/autofilter
# ranges from A1:D119 to A1:D118, and drops the now-unused shared strings.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 109 is "Worksheet | getRange | excel-range-areas | colorAllFormulaCells"
$ws.Rows.Item(109).Delete()

# Reflect the author's final cursor position/selection in the sheet.
$ws.Range("O91").Select()
